$wb = $excel.ActiveWorkbook

# Update the trade category value on the "licenseClosure" sheet (cell C2)
# from "Veterinary Trades" to "Flammables" as part of adding renewal scenarios.
$wsClosure = $wb.Worksheets.Item("licenseClosure")
$wsClosure.Range("C2").Value = "Flammables"

# Update the active sheet / selection state to match the new workbook state:
# the previously active sheet (tradeDetails) loses its selection marker and
# its last selection moves to F27, while licenseClosure becomes the active
# sheet with selection C7.
$wsTradeDetails = $wb.Worksheets.Item("tradeDetails")
$wsTradeDetails.Range("F27").Select()

$wsClosure.Activate()
$wsClosure.Range("C7").Select()
